$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update Status text for row 3 ("Ready for handoff" -> "Handback transform failed")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Populate Error Detail column (P) for row 3 on zh-cn and de-de sheets
$wsZhCn.Range("P3").Value = "Handback file name: lykmdecf.3qa is different with handoff file name: 02aca772-7870-4529-bd7c-7f0d236c9e64.5d484b6c5a66a87486b6993e1e206a53a2a3a877.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: lykmdecf.3qa is different with handoff file name: 02aca772-7870-4529-bd7c-7f0d236c9e64.5d484b6c5a66a87486b6993e1e206a53a2a3a877.de-de."

# Widen the Error Detail column on both sheets to fit new content
# (COM ColumnWidth is in "characters"; Excel rounds to whole pixels, so
#  39.15 is the value that round-trips to a stored column width of 40)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.15
$wsDeDe.Columns.Item(16).ColumnWidth = 39.15
